# Updated cryptos list with GitHub Actions
# Applies the per-row Price (D) / Volume(1h) (E) refresh, plus the OKB/Stacks
# row swap (rows 39-40), to match the latest scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates whose new text is NOT number-like -----------
# These can be written with a plain .Value assignment; Excel will keep them
# as text on its own (extra "thousand-dot" grouping, or non-numeric glyphs).
$dPlain = @(
  ,("D2",  "58.404.17")
  ,("D3",  "2.491.41")
  ,("D9",  "2.511.48")
  ,("D14", "2.930.39")
  ,("D15", "58.359.35")
  ,("D18", "2.498.77")
  ,("D29", "0.0₃0757")
  ,("D51", "1.748.72")
)
foreach ($pair in $dPlain) {
  $ws.Range($pair[0]).Value = $pair[1]
}

# --- Price (column D) updates whose new text WOULD be auto-parsed as a ----
# number (dropping meaningful trailing zeros, e.g. "522.80" -> 522.8), so a
# leading apostrophe is used to force literal text, exactly as typing it by
# hand into Excel would require.
$dText = @(
  ,("D5",  "522.80")
  ,("D6",  "135.98")
  ,("D8",  "0.561")
  ,("D10", "0.0994")
  ,("D13", "0.342")
  ,("D16", "22.28")
  ,("D19", "10.76")
  ,("D20", "4.21")
  ,("D21", "322.75")
  ,("D22", "0.999")
  ,("D23", "5.79")
  ,("D24", "64.68")
  ,("D25", "0.416")
  ,("D26", "0.162")
  ,("D27", "0.995")
  ,("D28", "7.45")
  ,("D30", "169.76")
  ,("D31", "1.20")
  ,("D33", "1.71")
  ,("D36", "18.18")
  ,("D41", "0.811")
  ,("D42", "5.23")
  ,("D43", "279.14")
  ,("D44", "3.49")
  ,("D46", "124.98")
  ,("D47", "0.0913")
  ,("D48", "0.0495")
  ,("D50", "17.27")
)
foreach ($pair in $dText) {
  $ws.Range($pair[0]).Value = "'" + $pair[1]
}

# --- Volume(1h) (column E) updates -----------------------------------------
# All percentages carry surrounding padding/spaces and a trailing "%", so
# they are never mistaken for numbers and can always be set directly.
$eVals = @(
  ,("E2",  "  -0.80%  ")
  ,("E3",  "  -0.30%  ")
  ,("E4",  "  -0.03%  ")
  ,("E5",  "  -1.89%  ")
  ,("E6",  "  +0.73%  ")
  ,("E7",  "  -0.25%  ")
  ,("E8",  "  -1.19%  ")
  ,("E9",  "  +0.46%  ")
  ,("E10", "  -1.68%  ")
  ,("E11", "  -0.82%  ")
  ,("E12", "  -1.36%  ")
  ,("E13", "  -1.45%  ")
  ,("E14", "  -0.26%  ")
  ,("E15", "  -0.73%  ")
  ,("E16", "  -1.78%  ")
  ,("E17", "  -1.14%  ")
  ,("E18", "  +0.39%  ")
  ,("E19", "  -2.34%  ")
  ,("E20", "  -0.77%  ")
  ,("E21", "  +0.32%  ")
  ,("E22", "  -0.03%  ")
  ,("E23", "  -2.29%  ")
  ,("E24", "  -0.49%  ")
  ,("E25", "  -0.68%  ")
  ,("E26", "  -0.83%  ")
  ,("E27", "  -0.37%  ")
  ,("E28", "  -0.69%  ")
  ,("E29", "  -0.22%  ")
  ,("E30", "  -0.25%  ")
  ,("E31", "  +6.35%  ")
  ,("E32", "  -0.85%  ")
  ,("E33", "  -1.61%  ")
  ,("E35", "  -0.34%  ")
  ,("E36", "  -0.58%  ")
  ,("E37", "  -0.98%  ")
  ,("E38", "  +0.66%  ")
  ,("E39", "  -1.85%  ")
  ,("E40", "  -0.24%  ")
  ,("E41", "  +1.82%  ")
  ,("E42", "  +5.20%  ")
  ,("E43", "  -0.88%  ")
  ,("E44", "  -2.12%  ")
  ,("E45", "  +0.16%  ")
  ,("E46", "  -2.72%  ")
  ,("E47", "  -1.25%  ")
  ,("E48", "  -0.70%  ")
  ,("E49", "  -1.03%  ")
  ,("E50", "  +0.63%  ")
  ,("E51", "  +0.12%  ")
)
foreach ($pair in $eVals) {
  $ws.Range($pair[0]).Value = $pair[1]
}

# --- Rows 39 & 40 swapped places in the ranking (OKB <-> Stacks) ----------
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'1.48"

$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "'36.67"
